# Update column C ("Förändrad") date value from 2023-09-19 (45188) to
# 2023-09-20 (45189) for every data row on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 257 }

$range = $ws.Range("C2:C$lastRow")
for ($i = 1; $i -le $range.Rows.Count; $i++) {
    $cell = $range.Cells.Item($i, 1)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
